$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.273.29'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.667.13'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5281'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06367'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.91'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07832'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").Value = '1.669.37'
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").Value = '1.895.63'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5596'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.09%  '
$ws.Range("D16").Value = '0.0₅8111'
$ws.Range("E16").Value = '  -0.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.66'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '26.288.76'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.727'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '199.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.065'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.010'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1214'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.233'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.530'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05907'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.281'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.510'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.322'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.596'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9603'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.820'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5794'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.36%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("E40").Value = '  +1.15%  '
$ws.Range("D41").Value = '1.074.39'
$ws.Range("E41").Value = '  +2.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8585'
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.80'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("D45").Value = '1.805.57'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.014'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4412'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₈104'
$ws.Range("E49").Value = '  -3.55%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.083'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05144'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.39%  '
